# Update BMI PGS Calculation before sending to Nur.
#
# Inserts a new "N.of.non_zero_Variants" column right after "N.of.Variants"
# (the old Ambg / Total.Matched / Flipped.Matched columns shift one to the
# right), and refreshes the matching-summary figures for all five PGS
# traits with the recalculated numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new column at C; existing C..H shift to D..I.
$ws.Columns.Item(3).Insert()

# 2) Header + width for the newly inserted column.
$ws.Cells.Item(1, 3).Value = "N.of.non_zero_Variants"
$ws.Columns.Item(3).ColumnWidth = 21.8

# 3) Rewrite every data row (A:I) with the refreshed values.
#    A = Trait, B = N.of.Variants, C = N.of.non_zero_Variants,
#    D = N.of.Ambg, E = % Ambg, F = Total.Matched, G = % Total.Matched,
#    H = Flipped.Mathced, I = % Flipped.Matched

# Row 2 - PGS003844
$ws.Cells.Item(2, 1).Value = "PGS003844"
$ws.Cells.Item(2, 2).Value = 87
$ws.Cells.Item(2, 3).Value = 87
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 78
$ws.Cells.Item(2, 7).Value = 89.66
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0

# Row 3 - PGS003845
$ws.Cells.Item(3, 1).Value = "PGS003845"
$ws.Cells.Item(3, 2).Value = 752209
$ws.Cells.Item(3, 3).Value = 225012
$ws.Cells.Item(3, 4).Value = 334
$ws.Cells.Item(3, 5).Value = 0.04
$ws.Cells.Item(3, 6).Value = 215351
$ws.Cells.Item(3, 7).Value = 95.71
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0

# Row 4 - PGS003846
$ws.Cells.Item(4, 1).Value = "PGS003846"
$ws.Cells.Item(4, 2).Value = 817189
$ws.Cells.Item(4, 3).Value = 817189
$ws.Cells.Item(4, 4).Value = 67
$ws.Cells.Item(4, 5).Value = 0.01
$ws.Cells.Item(4, 6).Value = 771884
$ws.Cells.Item(4, 7).Value = 94.46
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0

# Row 5 - PGS003847
$ws.Cells.Item(5, 1).Value = "PGS003847"
$ws.Cells.Item(5, 2).Value = 300864
$ws.Cells.Item(5, 3).Value = 300864
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 293162
$ws.Cells.Item(5, 7).Value = 97.44
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0

# Row 6 - PGS003848
$ws.Cells.Item(6, 1).Value = "PGS003848"
$ws.Cells.Item(6, 2).Value = 741723
$ws.Cells.Item(6, 3).Value = 382148
$ws.Cells.Item(6, 4).Value = 79
$ws.Cells.Item(6, 5).Value = 0.01
$ws.Cells.Item(6, 6).Value = 361296
$ws.Cells.Item(6, 7).Value = 94.54
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(6, 9).Value = 0
